$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("score")

$names = @("Steven N", "Steven Neveadomi", "Tori", "Steven", "TORI IS MATLAB KING", "Shaleen", "Yo Dawg Crilla", "Steven N", "Steven", "NA")
$scores = @(740, 700, 300, 280, 260, 200, 170, 160, 110, 80)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $scores[$i]
}
